$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 347, shifting the existing rows 347:358 down to 348:359.
$ws.Rows.Item(347).Insert()

# Populate the newly inserted row 347 with the new record.
$ws.Cells.Item(347, 1).Value = 4
$ws.Cells.Item(347, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(347, 3).Value = "Los Lagos"
$ws.Cells.Item(347, 4).Value = 44747
$ws.Cells.Item(347, 5).Value = 10
$ws.Cells.Item(347, 6).Value = 100112023
$ws.Cells.Item(347, 7).Value = "Brócoli"
$ws.Cells.Item(347, 8).Value = "Sin especificar"
$ws.Cells.Item(347, 9).Value = "Primera"
$ws.Cells.Item(347, 10).Value = 1200
$ws.Cells.Item(347, 11).Value = 1500
$ws.Cells.Item(347, 12).Value = 1500
$ws.Cells.Item(347, 13).Value = 1500
$ws.Cells.Item(347, 14).Value = "$/unidad"
$ws.Cells.Item(347, 15).Value = "Región Metropolitana"
$ws.Cells.Item(347, 16).Value = 1500
$ws.Cells.Item(347, 17).Value = 1
$ws.Cells.Item(347, 18).Value = "Hortaliza"
